# Add a new "doSearch" worksheet after the existing "doLogIn" sheet,
# populate it with the search-term / watchlist values, make it the
# active sheet, and select A2 (mirrors a manual "test search and add to
# watchlist" recording in Excel).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after "doLogIn".
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "doSearch"

# Populate cells in the same order the strings were first typed so the
# shared-string table indices line up with the recorded session
# (Search Term, Stock, Etherium, Bitcoin, Bonds).
$ws2.Range("A1").Value = "Search Term "
$ws2.Range("A4").Value = "Stock"
$ws2.Range("A3").Value = "Etherium"
$ws2.Range("A2").Value = "Bitcoin "
$ws2.Range("A5").Value = "Bonds "

# Match the recorded column width for column A.
$ws2.Columns.Item(1).ColumnWidth = 12.83

# Leave A2 selected/active on the new sheet, and make it the active tab.
$ws2.Range("A2").Select() | Out-Null
$ws2.Activate() | Out-Null
